$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 45 new data rows (rows 102-146), continuing the existing
# regcntr_id / device_id / lang_code / is_active / cr_by / cr_dtimes pattern.
for ($i = 0; $i -lt 45; $i++) {
    $row = 102 + $i
    $ws.Cells.Item($row, 1).Value = 10002 + ($i % 9)
    $ws.Cells.Item($row, 2).Value = 3000121 + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Match the saved view: selection over the newly added rows.
[void]$ws.Range("A102:F146").Select()

# Portrait page orientation (as set in the source file's page setup).
$ws.PageSetup.Orientation = 1
